# Weekly data refresh: a new price observation is inserted at row 54
# (pushing the existing rows 54-121 down to 55-122), adding one more
# week of "Membrillo" (quince) price data for Feria Lagunitas de Puerto
# Montt. The sheet's used range grows from A1:T121 to A1:T122.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 54; Excel shifts rows
# 54:121 down to 55:122 and extends the sheet dimension automatically.
$ws.Rows("54:54").Insert()

# Populate the newly inserted row 54 with the new weekly record.
$ws.Cells.Item(54, 1).Value = 4
$ws.Cells.Item(54, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(54, 3).Value = "Los Lagos"
$ws.Cells.Item(54, 4).Value = 44994
$ws.Cells.Item(54, 5).Value = 10
$ws.Cells.Item(54, 6).Value = "Fruta"
$ws.Cells.Item(54, 7).Value = 100104
$ws.Cells.Item(54, 8).Value = "Frutos de pepita"
$ws.Cells.Item(54, 9).Value = 100104003
$ws.Cells.Item(54, 10).Value = "Membrillo"
$ws.Cells.Item(54, 11).Value = "Champion"
$ws.Cells.Item(54, 12).Value = "Primera"
$ws.Cells.Item(54, 13).Value = 400
$ws.Cells.Item(54, 14).Value = 17000
$ws.Cells.Item(54, 15).Value = 18000
$ws.Cells.Item(54, 16).Value = 17500
$ws.Cells.Item(54, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(54, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(54, 19).Value = 972
$ws.Cells.Item(54, 20).Value = 18

# Make sure the date cell keeps the same date/time number format as the
# other cells in column D (style index 2 in styles.xml).
$ws.Cells.Item(54, 4).NumberFormat = $ws.Cells.Item(55, 4).NumberFormat
